# table_1_78.xlsx edit: the "divide by 2" constant used for rows 35-38 in
# column D was changed to "divide by 3".
#
#   D35 (standalone formula)                : =$D$3/2  ->  =$D$3/3
#   D36 (master of shared formula D36:D38)  : =$D$3/2  ->  =$D$3/3
#
# D37/D38 are followers of D36's shared formula, so they pick up the new
# formula automatically. Columns E/K (which depend on D) and the grand
# total in C40 are plain formulas that recalculate from the new D values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D35 only ever held its own (non-shared) formula; update it on its own so
# it stays a standalone formula cell, matching how it was before.
$ws.Range("D35").Formula = '=$D$3/3'

# D36:D38 share one formula (D36 is the master cell). Assigning the same
# formula to the whole D36:D38 range at once preserves that shared-formula
# grouping (D37/D38 keep referencing D36's formula) instead of turning each
# cell into its own separate formula.
$ws.Range("D36:D38").Formula = '=$D$3/3'

# Force a full recalculation so every dependent cell (E35:E38, K35:K38,
# C40's FLOOR(SUM(K3:K38),1), etc.) is refreshed with the new values.
$excel.CalculateFullRebuild()

# Reflect the saved scroll position / selection from the edited workbook:
# the view had scrolled down and the active cell moved to E37.
$ws.Range("A25").Select()
$ws.Range("E37").Select()
